$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.155.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("E5").Value = "  -4.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5272"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.67%  "

$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06281"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07526"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.680.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.483"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5647"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008096"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.185.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.845"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.192"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.11%  "

$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1255"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.585"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.22%  "

$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06414"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.337"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.521"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.480"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.645"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.001"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6084"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.418"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.156"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.100.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01606"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8663"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.831.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05272"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.949"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4273"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("E51").Value = "  -2.42%  "

